{"js": "// Replace the three-digit \u00f7 one-digit division prompts in the table with\n// their new values, per the commit's regenerated numbers.\nconst replacements = [\n  [\"702\u00f76=\", \"441\u00f79=\"],\n  [\"332\u00f72=\", \"270\u00f72=\"],\n  [\"924\u00f74=\", \"768\u00f75=\"],\n  [\"957\u00f78=\", \"367\u00f77=\"],\n  [\"230\u00f76=\", \"742\u00f76=\"],\n  [\"529\u00f72=\", \"270\u00f75=\"],\n  [\"685\u00f74=\", \"572\u00f74=\"],\n  [\"352\u00f78=\", \"239\u00f78=\"],\n  [\"370\u00f74=\", \"325\u00f78=\"],\n  [\"238\u00f73=\", \"342\u00f79=\"],\n  [\"746\u00f72=\", \"783\u00f73=\"],\n  [\"257\u00f78=\", \"746\u00f74=\"],\n  [\"107\u00f76=\", \"174\u00f76=\"],\n  [\"365\u00f76=\", \"711\u00f73=\"],\n  [\"318\u00f77=\", \"532\u00f72=\"],\n  [\"985\u00f72=\", \"367\u00f73=\"],\n  [\"573\u00f74=\", \"233\u00f74=\"],\n  [\"464\u00f75=\", \"660\u00f78=\"],\n  [\"521\u00f78=\", \"851\u00f77=\"],\n  [\"634\u00f78=\", \"242\u00f76=\"],\n  [\"228\u00f76=\", \"218\u00f78=\"],\n  [\"354\u00f72=\", \"932\u00f77=\"],\n  [\"573\u00f78=\", \"929\u00f77=\"],\n  [\"399\u00f72=\", \"129\u00f74=\"],\n  [\"782\u00f78=\", \"446\u00f74=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text to replace: \"${oldText}\"`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the three-digit \u00f7 one-digit division prompts in the table with\n# their new values, per the commit's regenerated numbers.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"702\u00f76=\", \"441\u00f79=\"),\n    @(\"332\u00f72=\", \"270\u00f72=\"),\n    @(\"924\u00f74=\", \"768\u00f75=\"),\n    @(\"957\u00f78=\", \"367\u00f77=\"),\n    @(\"230\u00f76=\", \"742\u00f76=\"),\n    @(\"529\u00f72=\", \"270\u00f75=\"),\n    @(\"685\u00f74=\", \"572\u00f74=\"),\n    @(\"352\u00f78=\", \"239\u00f78=\"),\n    @(\"370\u00f74=\", \"325\u00f78=\"),\n    @(\"238\u00f73=\", \"342\u00f79=\"),\n    @(\"746\u00f72=\", \"783\u00f73=\"),\n    @(\"257\u00f78=\", \"746\u00f74=\"),\n    @(\"107\u00f76=\", \"174\u00f76=\"),\n    @(\"365\u00f76=\", \"711\u00f73=\"),\n    @(\"318\u00f77=\", \"532\u00f72=\"),\n    @(\"985\u00f72=\", \"367\u00f73=\"),\n    @(\"573\u00f74=\", \"233\u00f74=\"),\n    @(\"464\u00f75=\", \"660\u00f78=\"),\n    @(\"521\u00f78=\", \"851\u00f77=\"),\n    @(\"634\u00f78=\", \"242\u00f76=\"),\n    @(\"228\u00f76=\", \"218\u00f78=\"),\n    @(\"354\u00f72=\", \"932\u00f77=\"),\n    @(\"573\u00f78=\", \"929\u00f77=\"),\n    @(\"399\u00f72=\", \"129\u00f74=\"),\n    @(\"782\u00f78=\", \"446\u00f74=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        Write-Output \"WARNING: not found -> $oldText\"\n    }\n}\n"}
